# Insert a new "is_targeted list" sheet (right after "assay_type list"),
# populate it with the boolean options as literal text "TRUE"/"FALSE",
# and repoint the N2:N1048576 data validation on the main sheet to use
# that new list instead of the inline "TRUE,FALSE" formula.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("assay_type list")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "is_targeted list"

# Force text (not boolean) values in the new list sheet.
$newSheet.Range("A1:A2").NumberFormat = "@"
$newSheet.Range("A1").Value = "'TRUE"
$newSheet.Range("A2").Value = "'FALSE"

# Point the is_targeted column's validation at the new list sheet instead
# of the inline "TRUE,FALSE" formula.
$mainSheet = $wb.Worksheets.Item("Export as TSV")
$targetRange = $mainSheet.Range("N2:N1048576")
$validation = $targetRange.Validation
$validation.Modify(3, 1, 1, "'is_targeted list'!`$A`$1:`$A`$2")
$validation.ErrorTitle = "Value must come from list"
$validation.ErrorMessage = "Value must be one of: TRUE / FALSE."
